$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 55

# Set the values for the new row first
$ws.Range("A$row").Value = 54
$ws.Range("B$row").Value = "denmark"
$ws.Range("C$row").Value = "superliga"
$ws.Range("D$row").Value = "2023-2024"
$ws.Range("E$row").Value = 45194.79166666666
$ws.Range("F$row").Value = "Nordsjaelland"
$ws.Range("G$row").Value = 0
$ws.Range("H$row").Value = "Hvidovre IF"
$ws.Range("I$row").Value = 0
$ws.Range("J$row").Value = 1.23
$ws.Range("K$row").Value = "18/09/2023 07:42"
$ws.Range("L$row").Value = 1.27
$ws.Range("M$row").Value = "25/09/2023 18:57"
$ws.Range("N$row").Value = 6.62
$ws.Range("O$row").Value = "18/09/2023 07:42"
$ws.Range("P$row").Value = 6.3
$ws.Range("Q$row").Value = "25/09/2023 18:59"
$ws.Range("R$row").Value = 11.54
$ws.Range("S$row").Value = "18/09/2023 07:42"
$ws.Range("T$row").Value = 9.949999999999999
$ws.Range("U$row").Value = "25/09/2023 18:59"
$ws.Range("V$row").Value = "https://www.betexplorer.com/football/denmark/superliga/nordsjaelland-hvidovre-if/dMp9qJ10/"

# Match the formatting used by the previous data row (row 54):
# column A is bold / bordered / centered, column E uses the custom
# date-time number format. Copy formats only, so no new style entries
# are introduced and the existing cellXfs are reused.
$ws.Range("A54").Copy()
$ws.Range("A$row").PasteSpecial(-4122)

$ws.Range("E54").Copy()
$ws.Range("E$row").PasteSpecial(-4122)

$excel.CutCopyMode = 0
